# Update the text in cell A1 (shared string "test" -> "Text in cell A1")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "Text in cell A1"

# Set column A width (diff shows XML width="14.85546875", i.e. ~104px /
# 14.14 "characters" at the standard Calibri 11 digit width). The
# ColumnWidth COM property expects a width in characters; 14 is the
# closest value that reproduces the custom width stored for column A.
$ws.Columns.Item(1).ColumnWidth = 14
